$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{Cell="D2"; Value="60.476.68"},
    @{Cell="E2"; Value="  -2.36%  "},
    @{Cell="D3"; Value="2.577.73"},
    @{Cell="E3"; Value="  -4.01%  "},
    @{Cell="E4"; Value="  -0.06%  "},
    @{Cell="D5"; Value="507.25"},
    @{Cell="E5"; Value="  -2.08%  "},
    @{Cell="D6"; Value="154.94"},
    @{Cell="E6"; Value="  -3.85%  "},
    @{Cell="D7"; Value="0.998"},
    @{Cell="E7"; Value="  -0.25%  "},
    @{Cell="D8"; Value="0.581"},
    @{Cell="E8"; Value="  -5.56%  "},
    @{Cell="D9"; Value="2.586.09"},
    @{Cell="E9"; Value="  -4.01%  "},
    @{Cell="D10"; Value="6.57"},
    @{Cell="E10"; Value="  +6.59%  "},
    @{Cell="E11"; Value="  -3.01%  "},
    @{Cell="E12"; Value="  -1.60%  "},
    @{Cell="E13"; Value="  +0.82%  "},
    @{Cell="D14"; Value="3.029.13"},
    @{Cell="E14"; Value="  -4.26%  "},
    @{Cell="D15"; Value="60.471.09"},
    @{Cell="E15"; Value="  -1.67%  "},
    @{Cell="D16"; Value="21.65"},
    @{Cell="E16"; Value="  -4.39%  "},
    @{Cell="E17"; Value="  -2.06%  "},
    @{Cell="D18"; Value="2.583.14"},
    @{Cell="E18"; Value="  -4.22%  "},
    @{Cell="D19"; Value="4.76"},
    @{Cell="E19"; Value="  -1.76%  "},
    @{Cell="D20"; Value="345.63"},
    @{Cell="E20"; Value="  -3.86%  "},
    @{Cell="D21"; Value="10.48"},
    @{Cell="E21"; Value="  -2.06%  "},
    @{Cell="D22"; Value="6.11"},
    @{Cell="E22"; Value="  -2.78%  "},
    @{Cell="D23"; Value="0.999"},
    @{Cell="E23"; Value="  -0.09%  "},
    @{Cell="D24"; Value="60.26"},
    @{Cell="E24"; Value="  -1.52%  "},
    @{Cell="D25"; Value="0.420"},
    @{Cell="E25"; Value="  -2.61%  "},
    @{Cell="E26"; Value="  -2.30%  "},
    @{Cell="D27"; Value="2.697.19"},
    @{Cell="E27"; Value="  -3.33%  "},
    @{Cell="D28"; Value="0.998"},
    @{Cell="E28"; Value="  -0.27%  "},
    @{Cell="D29"; Value="0.0₃0845"},
    @{Cell="E29"; Value="  -2.88%  "},
    @{Cell="D30"; Value="7.38"},
    @{Cell="E30"; Value="  -3.42%  "},
    @{Cell="E31"; Value="  -0.10%  "},
    @{Cell="D32"; Value="19.38"},
    @{Cell="E32"; Value="  -2.46%  "},
    @{Cell="D33"; Value="152.99"},
    @{Cell="E33"; Value="  -3.51%  "},
    @{Cell="E34"; Value="  -2.83%  "},
    @{Cell="D35"; Value="5.71"},
    @{Cell="E35"; Value="  +0.18%  "},
    @{Cell="D36"; Value="4.01"},
    @{Cell="E36"; Value="  -1.10%  "},
    @{Cell="E37"; Value="  -3.75%  "},
    @{Cell="E38"; Value="  +0.50%  "},
    @{Cell="E39"; Value="  -1.72%  "},
    @{Cell="D40"; Value="0.845"},
    @{Cell="E40"; Value="  -3.89%  "},
    @{Cell="D41"; Value="36.18"},
    @{Cell="E41"; Value="  +0.39%  "},
    @{Cell="E42"; Value="  -1.54%  "},
    @{Cell="D43"; Value="297.02"},
    @{Cell="E43"; Value="  -2.98%  "},
    @{Cell="D44"; Value="0.621"},
    @{Cell="E44"; Value="  -4.16%  "},
    @{Cell="D45"; Value="0.0999"},
    @{Cell="E45"; Value="  -2.46%  "},
    @{Cell="E46"; Value="  -0.27%  "},
    @{Cell="D47"; Value="0.0557"},
    @{Cell="E47"; Value="  -5.22%  "},
    @{Cell="D48"; Value="19.76"},
    @{Cell="E48"; Value="  -3.68%  "},
    @{Cell="E49"; Value="  -3.73%  "},
    @{Cell="D50"; Value="0.0233"},
    @{Cell="E50"; Value="  -3.09%  "},
    @{Cell="E51"; Value="  +0.02%  "}
)

foreach ($chg in $changes) {
    $cell = $ws.Range($chg.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $chg.Value
    $cell.Style = "Normal"
}
